$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (M1_PH)
$ws.Range("B2").Value = 1.255084099591062
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = -0.07526281778284652
$ws.Range("G2").Value = 0.06112493876441671

# Row 3 (CM2_PH)
$ws.Range("B3").Value = -0.0
$ws.Range("F3").Value = 0.1146013691737822

# Row 4 (CMN3_PH)
$ws.Range("F4").Value = -0.1651586729926508

# Row 5 (CMN4_PH)
$v = [double]"1.070040597132289e-15"
$ws.Range("B5").Value = $v
$ws.Range("C5").Value = -0.0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 2.343275768088881
$ws.Range("G5").Value = 0
